$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '34.943.10'
$ws.Cells.Item(3, 4).Value = '1.820.00'
$ws.Cells.Item(4, 5).Value = '  +0.30%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = '@'
$cell.Value = '231.47'
$cell.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +3.08%  '
$ws.Cells.Item(6, 5).Value = '  +1.66%  '
$ws.Cells.Item(7, 5).Value = '  +0.25%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = '@'
$cell.Value = '41.98'
$cell.Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +1.72%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.309'
$cell.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +6.70%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0684'
$cell.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +2.77%  '
$ws.Cells.Item(11, 5).Value = '  +0.56%  '
$ws.Cells.Item(12, 4).Value = '2.084.68'
$ws.Cells.Item(12, 5).Value = '  +1.26%  '
$ws.Cells.Item(13, 4).Value = '1.815.46'
$ws.Cells.Item(13, 5).Value = '  +1.14%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = '@'
$cell.Value = '11.07'
$cell.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +2.36%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.659'
$cell.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +5.82%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.65'
$cell.Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +6.25%  '
$ws.Cells.Item(17, 4).Value = '34.919.43'
$ws.Cells.Item(17, 5).Value = '  +1.53%  '
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = '@'
$cell.Value = '69.48'
$cell.Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +3.30%  '
$ws.Cells.Item(19, 4).Value = '0.0₃0786'
$ws.Cells.Item(19, 5).Value = '  +2.92%  '
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = '@'
$cell.Value = '238.08'
$cell.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -0.61%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = '@'
$cell.Value = '11.76'
$cell.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +6.20%  '
$ws.Cells.Item(22, 2).Value = 'Uniswap'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.60'
$cell.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +12.50%  '
$ws.Cells.Item(23, 2).Value = 'Dai'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.01'
$cell.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +0.34%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.24'
$cell.Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +4.41%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = '@'
$cell.Value = '172.34'
$cell.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +0.54%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = '@'
$cell.Value = '7.76'
$cell.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +1.86%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = '@'
$cell.Value = '17.41'
$cell.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +0.63%  '
$ws.Cells.Item(28, 5).Value = '  +0.19%  '
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.54'
$cell.Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +25.76%  '
$ws.Cells.Item(30, 5).Value = '  +0.32%  '
$ws.Cells.Item(31, 4).Value = '3.341.63'
$ws.Cells.Item(31, 5).Value = '  +37.53%  '
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0550'
$cell.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +7.81%  '
$ws.Cells.Item(33, 5).Value = '  +2.77%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.98'
$cell.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +3.76%  '
$ws.Cells.Item(35, 5).Value = '  +1.23%  '
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = '@'
$cell.Value = '92.59'
$cell.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +8.47%  '
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.675'
$cell.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +5.35%  '
$ws.Cells.Item(38, 5).Value = '  +5.24%  '
$ws.Cells.Item(39, 4).Value = '1.313.55'
$ws.Cells.Item(39, 5).Value = '  -0.25%  '
$ws.Cells.Item(40, 5).Value = '  +2.25%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.27'
$cell.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +2.88%  '
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.985'
$cell.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +5.25%  '
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = '@'
$cell.Value = '14.66'
$cell.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +0.00%  '
$ws.Cells.Item(44, 2).Value = 'RenderToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.29'
$cell.Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -1.89%  '
$ws.Cells.Item(45, 2).Value = 'HuobiToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.44'
$cell.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +0.33%  '
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.76'
$cell.Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -0.77%  '
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.17'
$cell.Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +5.77%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0510'
$cell.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -1.69%  '
$ws.Cells.Item(49, 4).Value = '1.996.71'
$ws.Cells.Item(49, 5).Value = '  +1.93%  '
$ws.Cells.Item(50, 5).Value = '  +0.21%  '
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = '@'
$cell.Value = '100.33'
$cell.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -0.26%  '
